$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell $ws 'D2' '22.427.73'
Set-TextCell $ws 'D3' '1.566.57'
Set-TextCell $ws 'E3' '  -0.38%  '
Set-TextCell $ws 'E4' '  +0.03%  '
Set-TextCell $ws 'E5' '  +0.04%  '
Set-TextCell $ws 'D6' '285.01'
Set-TextCell $ws 'E6' '  -2.34%  '
Set-TextCell $ws 'D7' '0.3630'
Set-TextCell $ws 'E7' '  -2.52%  '
Set-TextCell $ws 'D8' '48.49'
Set-TextCell $ws 'E8' '  -2.96%  '
Set-TextCell $ws 'D9' '0.3319'
Set-TextCell $ws 'E9' '  -2.06%  '
Set-TextCell $ws 'E10' '  -2.08%  '
Set-TextCell $ws 'D11' '0.07380'
Set-TextCell $ws 'E11' '  -2.48%  '
Set-TextCell $ws 'D12' '1.003'
Set-TextCell $ws 'E12' '  +0.09%  '
Set-TextCell $ws 'E13' '  -2.45%  '
Set-TextCell $ws 'E14' '  -1.35%  '
Set-TextCell $ws 'D15' '6.890'
Set-TextCell $ws 'E15' '  -0.92%  '
Set-TextCell $ws 'D16' '1.566.39'
Set-TextCell $ws 'E16' '  -0.28%  '
Set-TextCell $ws 'E17' '  -1.75%  '
Set-TextCell $ws 'D18' '87.89'
Set-TextCell $ws 'E18' '  -3.45%  '
Set-TextCell $ws 'D19' '0.06715'
Set-TextCell $ws 'E19' '  -0.55%  '
Set-TextCell $ws 'E20' '  +0.03%  '
Set-TextCell $ws 'D21' '6.328'
Set-TextCell $ws 'E21' '  +0.38%  '
Set-TextCell $ws 'D22' '16.19'
Set-TextCell $ws 'E22' '  -0.60%  '
Set-TextCell $ws 'D23' '11.98'
Set-TextCell $ws 'E23' '  -1.48%  '
Set-TextCell $ws 'D24' '22.411.82'
Set-TextCell $ws 'E24' '  -0.11%  '
Set-TextCell $ws 'D25' '2.377'
Set-TextCell $ws 'E25' '  +1.85%  '
Set-TextCell $ws 'D26' '2.533'
Set-TextCell $ws 'E26' '  -5.57%  '
Set-TextCell $ws 'D27' '150.38'
Set-TextCell $ws 'E27' '  +1.19%  '
Set-TextCell $ws 'D28' '19.39'
Set-TextCell $ws 'E28' '  -3.41%  '
Set-TextCell $ws 'D29' '4.998'
Set-TextCell $ws 'E29' '  -0.67%  '
Set-TextCell $ws 'D30' '123.81'
Set-TextCell $ws 'E30' '  -1.38%  '
Set-TextCell $ws 'D31' '1.743.00'
Set-TextCell $ws 'E31' '  -0.31%  '
Set-TextCell $ws 'D32' '1.034'
Set-TextCell $ws 'E32' '  -1.99%  '
Set-TextCell $ws 'D33' '2.007'
Set-TextCell $ws 'E33' '  +1.05%  '
Set-TextCell $ws 'E34' '  -1.38%  '
Set-TextCell $ws 'D35' '9.713'
Set-TextCell $ws 'E35' '  -1.13%  '
Set-TextCell $ws 'D36' '0.08240'
Set-TextCell $ws 'E36' '  -1.52%  '
Set-TextCell $ws 'D37' '0.02408'
Set-TextCell $ws 'E37' '  -2.88%  '
Set-TextCell $ws 'B38' 'Hedera'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D38' '0.06440'
Set-TextCell $ws 'E38' '  -1.06%  '
Set-TextCell $ws 'B39' 'Algorand'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws 'D39' '0.2236'
Set-TextCell $ws 'E39' '  -3.08%  '
Set-TextCell $ws 'D40' '5.361'
Set-TextCell $ws 'E40' '  -1.84%  '
Set-TextCell $ws 'E41' '  -5.10%  '
Set-TextCell $ws 'D42' '0.6252'
Set-TextCell $ws 'E42' '  +0.49%  '
Set-TextCell $ws 'E43' '  -1.51%  '
Set-TextCell $ws 'D44' '1.001'
Set-TextCell $ws 'E44' '  -0.01%  '
Set-TextCell $ws 'D45' '13.74'
Set-TextCell $ws 'E45' '  -1.99%  '
Set-TextCell $ws 'D46' '0.6059'
Set-TextCell $ws 'E46' '  +4.26%  '
Set-TextCell $ws 'D47' '3.748'
Set-TextCell $ws 'E47' '  -1.63%  '
Set-TextCell $ws 'D48' '2.025'
Set-TextCell $ws 'E48' '  -2.02%  '
Set-TextCell $ws 'D49' '123.27'
Set-TextCell $ws 'E49' '  -5.36%  '
Set-TextCell $ws 'D50' '1.211'
Set-TextCell $ws 'E50' '  -0.77%  '
Set-TextCell $ws 'D51' '0.07202'
Set-TextCell $ws 'E51' '  -1.55%  '
